# Daily attendance processing - 2025-12-14 19:00:15
#
# Re-order the "Recorded By" (column G) contributor lists so that real
# reviewer identities (backup@backdoor.com / dnasr281@gmail.com) are
# listed first, ahead of the generic "System"/"system" markers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tokens that should be promoted to the front of the comma-separated list,
# preserving their relative order. Every other token (System, system,
# admin@admin.com, ...) keeps its relative order and is appended after.
$priorityTokens = @("backup@backdoor.com", "dnasr281@gmail.com")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $value = $cell.Value()

    if ([string]::IsNullOrEmpty($value)) {
        continue
    }

    $tokens = $value -split ", "

    $front = @()
    $back = @()
    foreach ($tok in $tokens) {
        if ($priorityTokens -contains $tok) {
            $front += $tok
        } else {
            $back += $tok
        }
    }

    $newTokens = $front + $back
    $newValue = [string]::Join(", ", $newTokens)

    if ($newValue -ne $value) {
        $cell.Value = $newValue
    }
}
